$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-22 06:40:16"

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
